$d = $word.ActiveDocument

$p = $d.Paragraphs.Last
$start = $p.Range.Start
$end = $p.Range.End - 1

$r = $d.Range($start, $end)
$r.Text = ""

$r2 = $d.Range($start, $start)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00290245" w:rsidRPr="00371BEB" w:rsidRDefault="00290245" w:rsidP="00371BEB"><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-BE"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-BE"/></w:rPr><w:t xml:space="preserve">Kevin + </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-BE"/></w:rPr><w:t>Alessio</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/><w:lang w:val="fr-BE"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/><w:lang w:val="fr-BE"/></w:rPr><w:t>2. Présentation du TFE</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/><w:lang w:val="fr-BE"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/><w:lang w:val="fr-BE"/></w:rPr><w:t>2.1 Avant</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-BE"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-BE"/></w:rPr><w:t xml:space="preserve">Tout commence lorsqu’un client de Technord introduit une demande pour la réalisation d’une application aidant l’ouverture de chantier. En effet, celles-ci peuvent devenir facilement compliquées à faire proprement au vu du grand nombre de variables à prendre en compte. Il faut </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="yellow"/><w:lang w:val="fr-BE"/></w:rPr><w:t>keep</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="yellow"/><w:lang w:val="fr-BE"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="yellow"/><w:lang w:val="fr-BE"/></w:rPr><w:t>track</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-BE"/></w:rPr><w:t xml:space="preserve"> de chaque membre du personnel interne de l’entreprise participant à l’ouverture de chantier, des habilitations de ceux-ci, des sociétés sous-traitantes, des certifications de ces dernières et ainsi de suite.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-BE"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-BE"/></w:rPr><w:t xml:space="preserve">Pour répondre à ce besoin, l’équipe MES de chez Technord a développé une application en ASP.Net avec l’aide du framework Telerik qui s’occupera de faciliter tout le processus mentionner précédemment. Durant les années qui vont suivre, d’autres clients vont faire des demandes similaires, ce qui va conduire Technord a dès lors proposé l’application déjà existante à tous ceux le voulant, avec </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-BE"/></w:rPr><w:t>à chaque fois un peu de personnalisation ajoutée au projet. Ces clients vont aussi demander de nouvelles fonctionnalités, qui vont être petit à petit ajoutées au projet initial.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/><w:lang w:val="fr-BE"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/><w:lang w:val="fr-BE"/></w:rPr><w:t>2.2 La problématique</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-BE"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-BE"/></w:rPr><w:t xml:space="preserve">L’application a été développée sur plusieurs années, par différentes personnes et avec des technologies déjà relativement anciennes pour l’époque. De plus, n’étant pas toujours la priorité de l’entreprise, elle fut réalisée sans forcément porter d’attention à la bonne réalisation des choses, la rendant remplie de bug et peu attirante visuellement. Par ailleurs, </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-BE"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">le processus de création s’étant fait de </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="yellow"/><w:lang w:val="fr-BE"/></w:rPr><w:t>manière dispersée</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-BE"/></w:rPr><w:t>, le code de l’application est désordonné, rendant l’ajout de nouvelles fonctionnalités plus difficile que ce que ça ne devrait l’être. Pour toutes ces raisons, l’équipe de Technord a décidé de tout recommencer de zéro, avec de nouvelles technologies et une interface grandement améliorée</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/><w:lang w:val="fr-BE"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/><w:lang w:val="fr-BE"/></w:rPr><w:t>2.3 Le projet</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-BE"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-BE"/></w:rPr><w:t>Ce qui nous a été demandé, à Noah et moi, est de migrer l’entièreté de l’ancienne solution monolithe vers une application composé de deux services comme l’est maintenant courant : un front-end et un back-end.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-BE"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-BE"/></w:rPr><w:t xml:space="preserve">Le front-end permettra à l’utilisateur d’interagir avec les fonctionnalités de l’application et sera réalisé en </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-BE"/></w:rPr><w:t>Angular</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-BE"/></w:rPr><w:t xml:space="preserve"> à l’aide du framework Kendo</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-BE"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-BE"/></w:rPr><w:t>Le back-end servira d’intermédiaire entre l’interface visuelle et la base de données et s’occupera du traitement des requêtes. Il sera réalisé en .Net Framework avec le langage de programmation C#.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-BE"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-BE"/></w:rPr><w:t>La base de données restera majoritairement inchangée pour permettre une installation simple de la nouvelle solution chez le client.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-BE"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="yellow"/><w:lang w:val="fr-BE"/></w:rPr><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="yellow"/><w:lang w:val="fr-BE"/></w:rPr><w:t>metnionner</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="yellow"/><w:lang w:val="fr-BE"/></w:rPr><w:t xml:space="preserve"> le </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="yellow"/><w:lang w:val="fr-BE"/></w:rPr><w:t>figma</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="yellow"/><w:lang w:val="fr-BE"/></w:rPr><w:t> </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="yellow"/><w:lang w:val="fr-BE"/></w:rPr><w:t>,</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="yellow"/><w:lang w:val="fr-BE"/></w:rPr><w:t xml:space="preserve"> partie déroulement </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="yellow"/><w:lang w:val="fr-BE"/></w:rPr><w:t>?)</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/><w:lang w:val="fr-BE"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/><w:lang w:val="fr-BE"/></w:rPr><w:t>2.4 Les objectifs</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-BE"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-BE"/></w:rPr><w:t>D’abord, il est très important que la totalité des fonctionnalités de l’ancienne solution soit disponibles sur la nouvelle. Il est question fournir une version améliorée aux clients. L’application sera donc testée rigoureusement afin de ne pas créer une perte de productions chez ces derniers.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-BE"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-BE"/></w:rPr><w:t>Ensuite, comme déjà mentionné, l</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-BE"/></w:rPr><w:t>a nouvelle application sera un renouvellement visuel. La barre d’acceptation sera donc grandement relevée en ce qui concerne le design de l’interface.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-BE"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-BE"/></w:rPr><w:t>Enfin, il devra être possible de facilement ajoutée de nouvelles fonctionnalités au projet, que ça soit dans un future lointain ou non. Certains ajouts vont d’ailleurs être proposés comme tâches bonus aux stagiaires en fonction de l’avancement du stage.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-BE"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r2.InsertXML($xml)
